$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.185.59'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.333.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '517.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +1.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.328.34'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.49%  '
$ws.Range('E10').Value = '  +8.12%  '
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('E12').Value = '  +7.71%  '
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.98'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.724.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '56.288.43'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.63%  '
$ws.Range('E17').Value = '  +4.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.315.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.24%  '
$ws.Range('E19').Value = '  +2.43%  '
$ws.Range('E20').Value = '  +3.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '320.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.66'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.994'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('E26').Value = '  +5.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '171.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.71'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.19'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0728'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.69%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.996'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('E36').Value = '  +5.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.924'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('E38').Value = '  +7.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.99%  '
$ws.Range('E40').Value = '  +8.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.382'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '138.97'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +11.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.57'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '272.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.04'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0508'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0927'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.55%  '
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('E50').Value = '  +4.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.73%  '
